$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 146, pushing existing rows 146:240 down to 147:241.
$ws.Rows.Item(146).Insert()

# Populate the newly inserted row 146 with its data.
$ws.Range("A146").Value2 = 11
$ws.Range("B146").Value = "Vega Monumental Concepción"
$ws.Range("C146").Value = "Bíobío"
$ws.Range("D146").Value2 = 44777
$ws.Range("E146").Value2 = 8
$ws.Range("F146").Value2 = 100112045
$ws.Range("G146").Value = "Zapallo"
$ws.Range("H146").Value = "Camote"
$ws.Range("I146").Value = "1a (guarda)"
$ws.Range("J146").Value2 = 400
$ws.Range("K146").Value2 = 700
$ws.Range("L146").Value2 = 750
$ws.Range("M146").Value2 = 725
$ws.Range("N146").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O146").Value = "Región de O'Higgins"
$ws.Range("P146").Value2 = 725
$ws.Range("Q146").Value2 = 1
$ws.Range("R146").Value = "Hortaliza"
